$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = "System, backup@backdoor.com, system"
    4 = "System, backup@backdoor.com"
    5 = "System, backup@backdoor.com"
    7 = "System, admin@admin.com"
    8 = "System, backup@backdoor.com"
    11 = "dnasr281@gmail.com, System"
    17 = "dnasr281@gmail.com, System"
    28 = "System, backup@backdoor.com, system"
    30 = "System, backup@backdoor.com"
    31 = "System, backup@backdoor.com"
    33 = "System, admin@admin.com"
    34 = "System, backup@backdoor.com"
    37 = "dnasr281@gmail.com, System"
    43 = "dnasr281@gmail.com, System"
    54 = "System, backup@backdoor.com, system"
    56 = "System, backup@backdoor.com"
    57 = "System, backup@backdoor.com"
    59 = "System, admin@admin.com"
    60 = "System, backup@backdoor.com"
    63 = "dnasr281@gmail.com, System"
    69 = "dnasr281@gmail.com, System"
    80 = "System, backup@backdoor.com"
    81 = "System, backup@backdoor.com"
    82 = "System, backup@backdoor.com"
    87 = "dnasr281@gmail.com, admin@admin.com"
    93 = "dnasr281@gmail.com, System"
    94 = "dnasr281@gmail.com, System"
    96 = "dnasr281@gmail.com, System"
    106 = "System, backup@backdoor.com"
    107 = "System, backup@backdoor.com"
    108 = "System, backup@backdoor.com"
    113 = "dnasr281@gmail.com, admin@admin.com"
    119 = "dnasr281@gmail.com, System"
    120 = "dnasr281@gmail.com, System"
    122 = "dnasr281@gmail.com, System"
    132 = "System, backup@backdoor.com"
    133 = "System, backup@backdoor.com"
    134 = "System, backup@backdoor.com"
    139 = "dnasr281@gmail.com, admin@admin.com"
    145 = "dnasr281@gmail.com, System"
    146 = "dnasr281@gmail.com, System"
    148 = "dnasr281@gmail.com, System"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
